# Insert 4 new rows at 765-768 (weekly update: new week of Manzana prices)
# shifts old rows 765-777 down to 769-781, matching the rest of the diff automatically
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A765:T768").EntireRow.Insert()

# Row 765
$ws.Cells.Item(765, 1).Value = 5
$ws.Cells.Item(765, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(765, 3).Value = 'Maule'
$ws.Cells.Item(765, 4).Value = 44656
$ws.Cells.Item(765, 5).Value = 7
$ws.Cells.Item(765, 6).Value = 'Fruta'
$ws.Cells.Item(765, 7).Value = 100104
$ws.Cells.Item(765, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(765, 9).Value = 100104002
$ws.Cells.Item(765, 10).Value = 'Manzana'
$ws.Cells.Item(765, 11).Value = 'Granny Smith'
$ws.Cells.Item(765, 12).Value = 'Especial'
$ws.Cells.Item(765, 13).Value = 230
$ws.Cells.Item(765, 14).Value = 9000
$ws.Cells.Item(765, 15).Value = 9000
$ws.Cells.Item(765, 16).Value = 9000
$ws.Cells.Item(765, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(765, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(765, 19).Value = 600
$ws.Cells.Item(765, 20).Value = 15

# Row 766
$ws.Cells.Item(766, 1).Value = 5
$ws.Cells.Item(766, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(766, 3).Value = 'Maule'
$ws.Cells.Item(766, 4).Value = 44656
$ws.Cells.Item(766, 5).Value = 7
$ws.Cells.Item(766, 6).Value = 'Fruta'
$ws.Cells.Item(766, 7).Value = 100104
$ws.Cells.Item(766, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(766, 9).Value = 100104002
$ws.Cells.Item(766, 10).Value = 'Manzana'
$ws.Cells.Item(766, 11).Value = 'Granny Smith'
$ws.Cells.Item(766, 12).Value = 'Primera'
$ws.Cells.Item(766, 13).Value = 200
$ws.Cells.Item(766, 14).Value = 8000
$ws.Cells.Item(766, 15).Value = 8000
$ws.Cells.Item(766, 16).Value = 8000
$ws.Cells.Item(766, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(766, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(766, 19).Value = 533
$ws.Cells.Item(766, 20).Value = 15

# Row 767
$ws.Cells.Item(767, 1).Value = 5
$ws.Cells.Item(767, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(767, 3).Value = 'Maule'
$ws.Cells.Item(767, 4).Value = 44656
$ws.Cells.Item(767, 5).Value = 7
$ws.Cells.Item(767, 6).Value = 'Fruta'
$ws.Cells.Item(767, 7).Value = 100104
$ws.Cells.Item(767, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(767, 9).Value = 100104002
$ws.Cells.Item(767, 10).Value = 'Manzana'
$ws.Cells.Item(767, 11).Value = 'Royal Gala'
$ws.Cells.Item(767, 12).Value = 'Primera'
$ws.Cells.Item(767, 13).Value = 180
$ws.Cells.Item(767, 14).Value = 8000
$ws.Cells.Item(767, 15).Value = 8000
$ws.Cells.Item(767, 16).Value = 8000
$ws.Cells.Item(767, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(767, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(767, 19).Value = 533
$ws.Cells.Item(767, 20).Value = 15

# Row 768
$ws.Cells.Item(768, 1).Value = 5
$ws.Cells.Item(768, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(768, 3).Value = 'Maule'
$ws.Cells.Item(768, 4).Value = 44656
$ws.Cells.Item(768, 5).Value = 7
$ws.Cells.Item(768, 6).Value = 'Fruta'
$ws.Cells.Item(768, 7).Value = 100104
$ws.Cells.Item(768, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(768, 9).Value = 100104002
$ws.Cells.Item(768, 10).Value = 'Manzana'
$ws.Cells.Item(768, 11).Value = 'Scarlett'
$ws.Cells.Item(768, 12).Value = 'Especial'
$ws.Cells.Item(768, 13).Value = 370
$ws.Cells.Item(768, 14).Value = 10000
$ws.Cells.Item(768, 15).Value = 11000
$ws.Cells.Item(768, 16).Value = 10676
$ws.Cells.Item(768, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(768, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(768, 19).Value = 712
$ws.Cells.Item(768, 20).Value = 15

Write-Output "done"
